# Auto-generated Excel COM-interop edit script
# Applies cached-value corrections to the Leve profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 199
$ws.Cells.Item(9, 10).Value = 356
$ws.Cells.Item(9, 12).Value = 356
$ws.Cells.Item(9, 14).Value = -694
$ws.Cells.Item(18, 8).Value = 560.5
$ws.Cells.Item(18, 9).Value = 615
$ws.Cells.Item(18, 10).Value = 469.66666
$ws.Cells.Item(18, 11).Value = 615
$ws.Cells.Item(18, 12).Value = 469.66666
$ws.Cells.Item(18, 13).Value = -331
$ws.Cells.Item(18, 14).Value = -1037.66666
$ws.Cells.Item(40, 8).Value = 5696.1924
$ws.Cells.Item(40, 9).Value = 7121.273
$ws.Cells.Item(40, 10).Value = 4651.1333
$ws.Cells.Item(40, 11).Value = 7121.273
$ws.Cells.Item(40, 12).Value = 4651.1333
$ws.Cells.Item(40, 13).Value = -6946.273
$ws.Cells.Item(40, 14).Value = -5001.1333
$ws.Cells.Item(43, 8).Value = 3629.9
$ws.Cells.Item(43, 10).Value = 3899.8572
$ws.Cells.Item(43, 12).Value = 3899.8572
$ws.Cells.Item(43, 14).Value = -4037.8572
$ws.Cells.Item(100, 8).Value = 4908.4062
$ws.Cells.Item(100, 9).Value = 3360.476
$ws.Cells.Item(100, 11).Value = 3360.476
$ws.Cells.Item(100, 13).Value = -2819.476
$ws.Cells.Item(132, 8).Value = 43546.69
$ws.Cells.Item(132, 9).Value = 3327.1428
$ws.Cells.Item(132, 10).Value = 90469.5
$ws.Cells.Item(132, 11).Value = 9981.428400000001
$ws.Cells.Item(132, 12).Value = 271408.5
$ws.Cells.Item(132, 13).Value = -7451.428400000001
$ws.Cells.Item(132, 14).Value = -276468.5
$ws.Cells.Item(137, 8).Value = 8011803.5
$ws.Cells.Item(137, 9).Value = 22250550
$ws.Cells.Item(137, 10).Value = 2509.1875
$ws.Cells.Item(137, 11).Value = 66751650
$ws.Cells.Item(137, 12).Value = 7527.5625
$ws.Cells.Item(137, 13).Value = -66749100
$ws.Cells.Item(137, 14).Value = -12627.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 230.9
$ws.Cells.Item(4, 9).Value = 230.9
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 230.9
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -114.9
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 13833.333
$ws.Cells.Item(23, 10).Value = 13833.333
$ws.Cells.Item(23, 12).Value = 13833.333
$ws.Cells.Item(23, 14).Value = -14351.333
$ws.Cells.Item(74, 8).Value = 3199.8333
$ws.Cells.Item(74, 9).Value = 3650.5
$ws.Cells.Item(74, 10).Value = 1622.5
$ws.Cells.Item(74, 11).Value = 3650.5
$ws.Cells.Item(74, 12).Value = 1622.5
$ws.Cells.Item(74, 13).Value = -2776.5
$ws.Cells.Item(74, 14).Value = -3370.5
$ws.Cells.Item(77, 8).Value = 3199.8333
$ws.Cells.Item(77, 9).Value = 3650.5
$ws.Cells.Item(77, 10).Value = 1622.5
$ws.Cells.Item(77, 11).Value = 18252.5
$ws.Cells.Item(77, 12).Value = 8112.5
$ws.Cells.Item(77, 13).Value = -13884.5
$ws.Cells.Item(77, 14).Value = -16848.5
$ws.Cells.Item(102, 8).Value = 2199.4
$ws.Cells.Item(102, 9).Value = 2280.4814
$ws.Cells.Item(102, 10).Value = 1469.6666
$ws.Cells.Item(102, 11).Value = 2280.4814
$ws.Cells.Item(102, 12).Value = 1469.6666
$ws.Cells.Item(102, 13).Value = -658.4814000000001
$ws.Cells.Item(102, 14).Value = -4713.6666
$ws.Cells.Item(122, 8).Value = 3783.9756
$ws.Cells.Item(122, 9).Value = 3392.861
$ws.Cells.Item(122, 10).Value = 6600
$ws.Cells.Item(122, 11).Value = 10178.583
$ws.Cells.Item(122, 12).Value = 19800
$ws.Cells.Item(122, 13).Value = -7728.582999999999
$ws.Cells.Item(122, 14).Value = -24700

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 33491
$ws.Cells.Item(22, 9).Value = 33491
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 33491
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -33318
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 45347
$ws.Cells.Item(86, 9).Value = 63983.438
$ws.Cells.Item(86, 11).Value = 63983.438
$ws.Cells.Item(86, 13).Value = -62860.438
$ws.Cells.Item(89, 8).Value = 45347
$ws.Cells.Item(89, 9).Value = 63983.438
$ws.Cells.Item(89, 11).Value = 319917.19
$ws.Cells.Item(89, 13).Value = -314301.19
$ws.Cells.Item(99, 8).Value = 29049.588
$ws.Cells.Item(99, 9).Value = 38676.637
$ws.Cells.Item(99, 11).Value = 38676.637
$ws.Cells.Item(99, 13).Value = -37178.637
$ws.Cells.Item(105, 8).Value = 2113
$ws.Cells.Item(105, 9).Value = 2266.7812
$ws.Cells.Item(105, 11).Value = 2266.7812
$ws.Cells.Item(105, 13).Value = -519.7811999999999
$ws.Cells.Item(107, 8).Value = 4626.227
$ws.Cells.Item(107, 9).Value = 4385.467
$ws.Cells.Item(107, 10).Value = 5142.143
$ws.Cells.Item(107, 11).Value = 4385.467
$ws.Cells.Item(107, 12).Value = 5142.143
$ws.Cells.Item(107, 13).Value = -2465.467
$ws.Cells.Item(107, 14).Value = -8982.143
$ws.Cells.Item(117, 8).Value = 89999.5
$ws.Cells.Item(117, 10).Value = 89999.5
$ws.Cells.Item(117, 12).Value = 89999.5
$ws.Cells.Item(117, 14).Value = -99177.5
$ws.Cells.Item(134, 8).Value = 5802
$ws.Cells.Item(134, 9).Value = 2862
$ws.Cells.Item(134, 11).Value = 8586
$ws.Cells.Item(134, 13).Value = -6051

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 456614.47
$ws.Cells.Item(22, 9).Value = 2636.125
$ws.Cells.Item(22, 10).Value = 1667223.4
$ws.Cells.Item(22, 11).Value = 2636.125
$ws.Cells.Item(22, 12).Value = 1667223.4
$ws.Cells.Item(22, 13).Value = -2286.125
$ws.Cells.Item(22, 14).Value = -1667923.4
$ws.Cells.Item(31, 8).Value = 6266.037
$ws.Cells.Item(31, 9).Value = 8223.6875
$ws.Cells.Item(31, 10).Value = 3418.5454
$ws.Cells.Item(31, 11).Value = 8223.6875
$ws.Cells.Item(31, 12).Value = 3418.5454
$ws.Cells.Item(31, 13).Value = -7928.6875
$ws.Cells.Item(31, 14).Value = -4008.5454
$ws.Cells.Item(34, 8).Value = 6266.037
$ws.Cells.Item(34, 9).Value = 8223.6875
$ws.Cells.Item(34, 10).Value = 3418.5454
$ws.Cells.Item(34, 11).Value = 8223.6875
$ws.Cells.Item(34, 12).Value = 3418.5454
$ws.Cells.Item(34, 13).Value = -8021.6875
$ws.Cells.Item(34, 14).Value = -3822.5454
$ws.Cells.Item(122, 8).Value = 3398.0476
$ws.Cells.Item(122, 9).Value = 4499.364
$ws.Cells.Item(122, 10).Value = 2186.6
$ws.Cells.Item(122, 11).Value = 13498.092
$ws.Cells.Item(122, 12).Value = 6559.799999999999
$ws.Cells.Item(122, 13).Value = -11048.092
$ws.Cells.Item(122, 14).Value = -11459.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 56333828
$ws.Cells.Item(4, 9).Value = 30235672
$ws.Cells.Item(4, 11).Value = 90707016
$ws.Cells.Item(4, 13).Value = -90706904
$ws.Cells.Item(5, 8).Value = 1517.5883
$ws.Cells.Item(5, 9).Value = 1711
$ws.Cells.Item(5, 11).Value = 5133
$ws.Cells.Item(5, 13).Value = -5021
$ws.Cells.Item(23, 8).Value = 136.76923
$ws.Cells.Item(23, 9).Value = 117.85714
$ws.Cells.Item(23, 10).Value = 158.83333
$ws.Cells.Item(23, 11).Value = 353.57142
$ws.Cells.Item(23, 12).Value = 476.49999
$ws.Cells.Item(23, 13).Value = -118.57142
$ws.Cells.Item(23, 14).Value = -946.49999
$ws.Cells.Item(33, 8).Value = 318.75
$ws.Cells.Item(33, 9).Value = 150
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 900
$ws.Cells.Item(33, 12).Value = 3600
$ws.Cells.Item(33, 13).Value = -617
$ws.Cells.Item(33, 14).Value = -4166
$ws.Cells.Item(68, 8).Value = 1443.4
$ws.Cells.Item(68, 9).Value = 1268.5
$ws.Cells.Item(68, 11).Value = 3805.5
$ws.Cells.Item(68, 13).Value = -2994.5
$ws.Cells.Item(71, 8).Value = 1443.4
$ws.Cells.Item(71, 9).Value = 1268.5
$ws.Cells.Item(71, 11).Value = 11416.5
$ws.Cells.Item(71, 13).Value = -7360.5
$ws.Cells.Item(113, 8).Value = 811.25
$ws.Cells.Item(113, 9).Value = 1182.5
$ws.Cells.Item(113, 10).Value = 625.625
$ws.Cells.Item(113, 11).Value = 3547.5
$ws.Cells.Item(113, 12).Value = 1876.875
$ws.Cells.Item(113, 13).Value = -1377.5
$ws.Cells.Item(113, 14).Value = -6216.875
$ws.Cells.Item(135, 8).Value = 1517.5883
$ws.Cells.Item(135, 9).Value = 1711
$ws.Cells.Item(135, 11).Value = 15399
$ws.Cells.Item(135, 13).Value = -12864

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 254.6923
$ws.Cells.Item(2, 9).Value = 146.11111
$ws.Cells.Item(2, 10).Value = 499
$ws.Cells.Item(2, 11).Value = 146.11111
$ws.Cells.Item(2, 12).Value = 499
$ws.Cells.Item(2, 13).Value = -33.11111
$ws.Cells.Item(2, 14).Value = -725
$ws.Cells.Item(95, 8).Value = 25288.834
$ws.Cells.Item(95, 10).Value = 25288.834
$ws.Cells.Item(95, 12).Value = 25288.834
$ws.Cells.Item(95, 14).Value = -30780.834
$ws.Cells.Item(97, 8).Value = 2866.6191
$ws.Cells.Item(97, 9).Value = 497.3421
$ws.Cells.Item(97, 10).Value = 25374.75
$ws.Cells.Item(97, 11).Value = 497.3421
$ws.Cells.Item(97, 12).Value = 25374.75
$ws.Cells.Item(97, 13).Value = -1.342100000000016
$ws.Cells.Item(97, 14).Value = -26366.75
$ws.Cells.Item(132, 8).Value = 5860.431
$ws.Cells.Item(132, 9).Value = 5190.7954
$ws.Cells.Item(132, 11).Value = 15572.3862
$ws.Cells.Item(132, 13).Value = -13042.3862

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1007.5
$ws.Cells.Item(55, 9).Value = 509.375
$ws.Cells.Item(55, 11).Value = 509.375
$ws.Cells.Item(55, 13).Value = -336.375
